$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update id and username/email, keep rest of the address data as-is
$ws.Range("A3").Value = 52
$ws.Range("B3").Value = "tejst5"
$ws.Range("C3").Value = "test5j@gmail.com"
$ws.Range("D3").Value = "Đạt"
$ws.Range("E3").Value = "Phan"
$ws.Range("F3").Value = 36572
$ws.Range("G3").Value = "Đường Ok"
$ws.Range("H3").Value = "Tam Thuấn"
$ws.Range("I3").Value = "Phúc Thọ"
$ws.Range("J3").Value = "Hà Nội"
$ws.Range("K3").Value = 2

# Row 4: becomes the previous "dat2k3" user, now with full address/dob filled in
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "dat2k3"
$ws.Range("C4").Value = "dat2k3@gmail.com"
$ws.Range("D4").Value = "Đạt"
$ws.Range("E4").Value = "Phan"
$ws.Range("F4").Value = 36572
$ws.Range("G4").Value = "Đường Ok"
$ws.Range("H4").Value = "Tam Thuấn"
$ws.Range("I4").Value = "Phúc Thọ"
$ws.Range("J4").Value = "Hà Nội"
$ws.Range("K4").Value = 2
